# Update inventory report: revise rows 2-4 content and remove rows 5-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to keep a text value (e.g. "3200000.00") instead of
    # being auto-coerced into a number, without leaving a residual
    # text-number-format style behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# Row 2: Laptop HP / Electrónica / 10 / 3200000.00 / Proveedor A
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Laptop HP"
$ws.Range("C2").Value = "Electrónica"
$ws.Range("D2").Value = 10
Set-TextValue $ws.Range("E2") "3200000.00"
$ws.Range("F2").Value = "Proveedor A"

# Row 3: Impresora Epson / Electrónica / 20 / 450000.00 / Proveedor B
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Impresora Epson"
$ws.Range("C3").Value = "Electrónica"
$ws.Range("D3").Value = 20
Set-TextValue $ws.Range("E3") "450000.00"
$ws.Range("F3").Value = "Proveedor B"

# Row 4: arroz / grano / 1 / 1750.00 / Proveedor B
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "arroz"
$ws.Range("C4").Value = "grano"
$ws.Range("D4").Value = 1
Set-TextValue $ws.Range("E4") "1750.00"
$ws.Range("F4").Value = "Proveedor B"

# Remove old rows 5, 6, 7 entirely (dimension shrinks from A1:F7 to A1:F4)
$ws.Range("A5:F7").Delete()
